# Auto-generated edit script applying numeric value updates
# to the Phoenix_Profits workbook, per the supplied OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 641.2
$ws.Range("I6").Value = 467
$ws.Range("K6").Value = 1401
$ws.Range("M6").Value = -1289
# Row 38
$ws.Range("H38").Value = 512.5
$ws.Range("J38").Value = 1916.6666
$ws.Range("L38").Value = 5749.9998
$ws.Range("N38").Value = -6493.9998
# Row 41
$ws.Range("H41").Value = 755.2941
$ws.Range("J41").Value = 414.57144
$ws.Range("L41").Value = 414.57144
$ws.Range("N41").Value = -1294.57144
# Row 48
$ws.Range("H48").Value = 9500
$ws.Range("J48").Value = 9500
$ws.Range("L48").Value = 28500
$ws.Range("N48").Value = -29084
# Row 56
$ws.Range("H56").Value = 9500
$ws.Range("J56").Value = 9500
$ws.Range("L56").Value = 28500
$ws.Range("N56").Value = -29568
# Row 58
$ws.Range("H58").Value = 818.1053000000001
$ws.Range("I58").Value = 73.5
$ws.Range("J58").Value = 1016.6667
$ws.Range("K58").Value = 220.5
$ws.Range("L58").Value = 3050.0001
$ws.Range("M58").Value = -70.5
$ws.Range("N58").Value = -3350.0001
# Row 107
$ws.Range("H107").Value = 4504999.5
$ws.Range("I107").Value = 4504999.5
$ws.Range("K107").Value = 4504999.5
$ws.Range("M107").Value = -4503079.5
# Row 138
$ws.Range("H138").Value = 2786.075
$ws.Range("I138").Value = 2088.7144
$ws.Range("J138").Value = 3161.577
$ws.Range("K138").Value = 6266.1432
$ws.Range("L138").Value = 9484.731
$ws.Range("M138").Value = -1126.1432
$ws.Range("N138").Value = -19764.731

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2997.5
$ws.Range("I2").Value = 2820.077
$ws.Range("J2").Value = 3766.3333
$ws.Range("K2").Value = 2820.077
$ws.Range("L2").Value = 3766.3333
$ws.Range("M2").Value = -2707.077
$ws.Range("N2").Value = -3992.3333
# Row 36
$ws.Range("H36").Value = 19000
$ws.Range("I36").Value = 19000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 19000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -18654
$ws.Range("N36").ClearContents()
# Row 61
$ws.Range("H61").Value = 3789.1143
$ws.Range("I61").Value = 2965.0908
$ws.Range("J61").Value = 5183.615
$ws.Range("K61").Value = 2965.0908
$ws.Range("L61").Value = 5183.615
$ws.Range("M61").Value = -2753.0908
$ws.Range("N61").Value = -5607.615
# Row 74
$ws.Range("H74").Value = 1452
$ws.Range("I74").Value = 1234.1666
$ws.Range("J74").Value = 2105.5
$ws.Range("K74").Value = 1234.1666
$ws.Range("L74").Value = 2105.5
$ws.Range("M74").Value = -360.1666
$ws.Range("N74").Value = -3853.5
# Row 77
$ws.Range("H77").Value = 1452
$ws.Range("I77").Value = 1234.1666
$ws.Range("J77").Value = 2105.5
$ws.Range("K77").Value = 6170.833000000001
$ws.Range("L77").Value = 10527.5
$ws.Range("M77").Value = -1802.833000000001
$ws.Range("N77").Value = -19263.5
# Row 97
$ws.Range("H97").Value = 1760.6364
$ws.Range("J97").Value = 3995.9092
$ws.Range("L97").Value = 3995.9092
$ws.Range("N97").Value = -4987.9092
# Row 116
$ws.Range("H116").Value = 2997.5
$ws.Range("I116").Value = 2820.077
$ws.Range("J116").Value = 3766.3333
$ws.Range("K116").Value = 2820.077
$ws.Range("L116").Value = 3766.3333
$ws.Range("M116").Value = -526.0770000000002
$ws.Range("N116").Value = -8354.3333
# Row 122
$ws.Range("H122").Value = 3288.7778
$ws.Range("I122").Value = 3288.7778
$ws.Range("K122").Value = 9866.3334
$ws.Range("M122").Value = -7416.3334
# Row 132
$ws.Range("H132").Value = 5283.2856
$ws.Range("I132").Value = 3744.375
$ws.Range("J132").Value = 7335.1665
$ws.Range("K132").Value = 11233.125
$ws.Range("L132").Value = 22005.4995
$ws.Range("M132").Value = -8703.125
$ws.Range("N132").Value = -27065.4995
# Row 136
$ws.Range("H136").Value = 3789.1143
$ws.Range("I136").Value = 2965.0908
$ws.Range("J136").Value = 5183.615
$ws.Range("K136").Value = 8895.2724
$ws.Range("L136").Value = 15550.845
$ws.Range("M136").Value = -6345.2724
$ws.Range("N136").Value = -20650.845

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2997.5
$ws.Range("I3").Value = 2820.077
$ws.Range("J3").Value = 3766.3333
$ws.Range("K3").Value = 2820.077
$ws.Range("L3").Value = 3766.3333
$ws.Range("M3").Value = -2706.077
$ws.Range("N3").Value = -3994.3333
# Row 80
$ws.Range("H80").Value = 339.5263
$ws.Range("J80").Value = 361.92856
$ws.Range("L80").Value = 361.92856
$ws.Range("N80").Value = -2357.92856
# Row 83
$ws.Range("H83").Value = 339.5263
$ws.Range("J83").Value = 361.92856
$ws.Range("L83").Value = 1809.6428
$ws.Range("N83").Value = -11793.6428

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3463.1072
$ws.Range("I31").Value = 1283.7778
$ws.Range("K31").Value = 1283.7778
$ws.Range("M31").Value = -988.7778000000001
# Row 34
$ws.Range("H34").Value = 3463.1072
$ws.Range("I34").Value = 1283.7778
$ws.Range("K34").Value = 1283.7778
$ws.Range("M34").Value = -1081.7778
# Row 60
$ws.Range("H60").Value = 38977
$ws.Range("J60").Value = 38977
$ws.Range("L60").Value = 38977
$ws.Range("N60").Value = -39999
# Row 62
$ws.Range("H62").Value = 3322.7334
$ws.Range("I62").Value = 3280.2307
$ws.Range("J62").Value = 3599
$ws.Range("K62").Value = 3280.2307
$ws.Range("L62").Value = 3599
$ws.Range("M62").Value = -2656.2307
$ws.Range("N62").Value = -4847
# Row 65
$ws.Range("H65").Value = 3322.7334
$ws.Range("I65").Value = 3280.2307
$ws.Range("J65").Value = 3599
$ws.Range("K65").Value = 16401.1535
$ws.Range("L65").Value = 17995
$ws.Range("M65").Value = -13281.1535
$ws.Range("N65").Value = -24235
# Row 94
$ws.Range("H94").Value = 1664.8667
$ws.Range("I94").Value = 1783.7
$ws.Range("K94").Value = 1783.7
$ws.Range("M94").Value = -1332.7

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 11243.286
$ws.Range("J55").Value = 15176.2
$ws.Range("L55").Value = 45528.60000000001
$ws.Range("N55").Value = -45882.60000000001
# Row 107
$ws.Range("H107").Value = 1555.2239
$ws.Range("J107").Value = 1690.0677
$ws.Range("L107").Value = 5070.203100000001
$ws.Range("N107").Value = -8910.203100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 967383.5
$ws.Range("I122").Value = 1190778.1
$ws.Range("J122").Value = 9978.143
$ws.Range("K122").Value = 3572334.3
$ws.Range("L122").Value = 29934.429
$ws.Range("M122").Value = -3569884.3
$ws.Range("N122").Value = -34834.429
# Row 126
$ws.Range("H126").Value = 2913.6
$ws.Range("J126").Value = 2433.125
$ws.Range("L126").Value = 7299.375
$ws.Range("N126").Value = -12239.375
# Row 132
$ws.Range("H132").Value = 7288.5557
$ws.Range("I132").Value = 8019.6
$ws.Range("J132").Value = 6374.75
$ws.Range("K132").Value = 24058.8
$ws.Range("L132").Value = 19124.25
$ws.Range("M132").Value = -21528.8
$ws.Range("N132").Value = -24184.25
# Row 135
$ws.Range("H135").Value = 95237.86
$ws.Range("J135").Value = 95237.86
$ws.Range("L135").Value = 95237.86
$ws.Range("N135").Value = -105377.86

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2298.7778
$ws.Range("I16").Value = 2298.7778
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2298.7778
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2128.7778
$ws.Range("N16").ClearContents()
# Row 40
$ws.Range("H40").Value = 4965.2915
$ws.Range("I40").Value = 4296.8
$ws.Range("K40").Value = 4296.8
$ws.Range("M40").Value = -4160.8
# Row 61
$ws.Range("H61").Value = 7580.3213
$ws.Range("I61").Value = 6997.5713
$ws.Range("J61").Value = 9328.571
$ws.Range("K61").Value = 6997.5713
$ws.Range("L61").Value = 9328.571
$ws.Range("M61").Value = -6795.5713
$ws.Range("N61").Value = -9732.571
# Row 113
$ws.Range("H113").Value = 7580.3213
$ws.Range("I113").Value = 6997.5713
$ws.Range("J113").Value = 9328.571
$ws.Range("K113").Value = 6997.5713
$ws.Range("L113").Value = 9328.571
$ws.Range("M113").Value = -4827.5713
$ws.Range("N113").Value = -13668.571
# Row 122
$ws.Range("H122").Value = 53334.332
$ws.Range("I122").Value = 53334.332
$ws.Range("K122").Value = 160002.996
$ws.Range("M122").Value = -157552.996
# Row 124
$ws.Range("H124").Value = 86653.28999999999
$ws.Range("J124").Value = 86653.28999999999
$ws.Range("L124").Value = 86653.28999999999
$ws.Range("N124").Value = -96473.28999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 128
$ws.Range("H128").Value = 69715
$ws.Range("J128").Value = 69715
$ws.Range("L128").Value = 69715
$ws.Range("N128").Value = -79675
# Row 132
$ws.Range("H132").Value = 4819.56
$ws.Range("I132").Value = 3023.9048
$ws.Range("K132").Value = 9071.714399999999
$ws.Range("M132").Value = -6541.714399999999
# Row 136
$ws.Range("H136").Value = 3742.5293
$ws.Range("I136").Value = 3939.5715
$ws.Range("K136").Value = 11818.7145
$ws.Range("M136").Value = -9268.7145

Write-Host "Updated $($wb.Worksheets.Count) sheets with new leve profit figures"
